# Edit: insert two new data rows (a new weekly "Primera"/"Segunda" pair for Brócoli)
# at the top of the existing data block, right before the current row 303.
# This pushes all subsequent rows down by 2 (old row 303 -> new row 305, etc.),
# and extends the used range from A1:R442 to A1:R444.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 303; Excel shifts rows 303:442 down to 305:444.
$ws.Rows("303:304").Insert()

# --- Populate new row 303 (quality "Primera") ---
$ws.Range("A303").Value = 8
$ws.Range("B303").Value = "Terminal La Palmera de La Serena"
$ws.Range("C303").Value = "Coquimbo"
$ws.Range("D303").Value = 44510
$ws.Range("E303").Value = 4
$ws.Range("F303").Value = 100112023
$ws.Range("G303").Value = "Brócoli"
$ws.Range("H303").Value = "Sin especificar"
$ws.Range("I303").Value = "Primera"
$ws.Range("J303").Value = 3000
$ws.Range("K303").Value = 600
$ws.Range("L303").Value = 700
$ws.Range("M303").Value = 650
$ws.Range("N303").Value = "`$/unidad"
$ws.Range("O303").Value = "Provincia del Elquí"
$ws.Range("P303").Value = 650
$ws.Range("Q303").Value = 1
$ws.Range("R303").Value = "Hortaliza"

# --- Populate new row 304 (quality "Segunda") ---
$ws.Range("A304").Value = 8
$ws.Range("B304").Value = "Terminal La Palmera de La Serena"
$ws.Range("C304").Value = "Coquimbo"
$ws.Range("D304").Value = 44510
$ws.Range("E304").Value = 4
$ws.Range("F304").Value = 100112023
$ws.Range("G304").Value = "Brócoli"
$ws.Range("H304").Value = "Sin especificar"
$ws.Range("I304").Value = "Segunda"
$ws.Range("J304").Value = 1600
$ws.Range("K304").Value = 500
$ws.Range("L304").Value = 550
$ws.Range("M304").Value = 525
$ws.Range("N304").Value = "`$/unidad"
$ws.Range("O304").Value = "Provincia del Elquí"
$ws.Range("P304").Value = 525
$ws.Range("Q304").Value = 1
$ws.Range("R304").Value = "Hortaliza"

# Keep date columns formatted the same way as the rest of column D.
$ws.Range("D303").NumberFormat = $ws.Range("D305").NumberFormat
$ws.Range("D304").NumberFormat = $ws.Range("D305").NumberFormat
